$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F8").Value = 7650
    $ws.Range("F9").Value = 71
    $ws.Range("F18").Value = 250
}
